# chore: update Sheets via scheduled runner
# Refreshes cached market-board price / profit figures on the Aegis_Profits
# workbook's per-job leve sheets (ALC, ARM, BSM, CUL, GSM, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 314567.6
$ws.Range("I33").Value = 45768.953
$ws.Range("J33").Value = 1300162.6
$ws.Range("K33").Value = 45768.953
$ws.Range("L33").Value = 1300162.6
$ws.Range("M33").Value = -45539.953
$ws.Range("N33").Value = -1300620.6

$ws.Range("H92").Value = 619.3913
$ws.Range("I92").Value = 605.8889
$ws.Range("K92").Value = 605.8889
$ws.Range("M92").Value = 642.1111

$ws.Range("H104").Value = 100000000
$ws.Range("I104").Value = 100000000
$ws.Range("K104").Value = 300000000
$ws.Range("M104").Value = -299998253

$ws.Range("H106").Value = 4668.4287
$ws.Range("I106").Value = 5133
$ws.Range("J106").Value = 4320
$ws.Range("K106").Value = 5133
$ws.Range("L106").Value = 4320
$ws.Range("M106").Value = -4502
$ws.Range("N106").Value = -5582

$ws.Range("H138").Value = 2079.8
$ws.Range("I138").Value = 1234.3889
$ws.Range("J138").Value = 2265.378
$ws.Range("K138").Value = 3703.1667
$ws.Range("L138").Value = 6796.134
$ws.Range("M138").Value = 1436.8333
$ws.Range("N138").Value = -17076.134

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35614.05
$ws.Range("I32").Value = 10651.437
$ws.Range("K32").Value = 10651.437
$ws.Range("M32").Value = -10364.437

$ws.Range("H74").Value = 792.4
$ws.Range("I74").Value = 587.3077
$ws.Range("J74").Value = 1173.2858
$ws.Range("K74").Value = 587.3077
$ws.Range("L74").Value = 1173.2858
$ws.Range("M74").Value = 286.6923
$ws.Range("N74").Value = -2921.2858

$ws.Range("H77").Value = 792.4
$ws.Range("I77").Value = 587.3077
$ws.Range("J77").Value = 1173.2858
$ws.Range("K77").Value = 2936.5385
$ws.Range("L77").Value = 5866.429
$ws.Range("M77").Value = 1431.4615
$ws.Range("N77").Value = -14602.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1144.5483
$ws.Range("I80").Value = 743.8182
$ws.Range("J80").Value = 1364.95
$ws.Range("K80").Value = 743.8182
$ws.Range("L80").Value = 1364.95
$ws.Range("M80").Value = 254.1818
$ws.Range("N80").Value = -3360.95

$ws.Range("H83").Value = 1144.5483
$ws.Range("I83").Value = 743.8182
$ws.Range("J83").Value = 1364.95
$ws.Range("K83").Value = 3719.091
$ws.Range("L83").Value = 6824.75
$ws.Range("M83").Value = 1272.909
$ws.Range("N83").Value = -16808.75

$ws.Range("H86").Value = 93892.164
$ws.Range("I86").Value = 124545.11
$ws.Range("J86").Value = 1933.3334
$ws.Range("K86").Value = 124545.11
$ws.Range("L86").Value = 1933.3334
$ws.Range("M86").Value = -123422.11
$ws.Range("N86").Value = -4179.3334

$ws.Range("H89").Value = 93892.164
$ws.Range("I89").Value = 124545.11
$ws.Range("J89").Value = 1933.3334
$ws.Range("K89").Value = 622725.55
$ws.Range("L89").Value = 9666.667
$ws.Range("M89").Value = -617109.55
$ws.Range("N89").Value = -20898.667

$ws.Range("H94").Value = 642.4286
$ws.Range("I94").Value = 574.5
$ws.Range("J94").Value = 733
$ws.Range("K94").Value = 574.5
$ws.Range("L94").Value = 733
$ws.Range("M94").Value = -123.5
$ws.Range("N94").Value = -1635

$ws.Range("H134").Value = 2313.5247
$ws.Range("I134").Value = 2221.5293
$ws.Range("J134").Value = 2782.7
$ws.Range("K134").Value = 6664.5879
$ws.Range("L134").Value = 8348.099999999999
$ws.Range("M134").Value = -4129.5879
$ws.Range("N134").Value = -13418.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1177.242
$ws.Range("I5").Value = 841.8571
$ws.Range("J5").Value = 1349.0244
$ws.Range("K5").Value = 2525.5713
$ws.Range("L5").Value = 4047.0732
$ws.Range("M5").Value = -2413.5713
$ws.Range("N5").Value = -4271.0732

$ws.Range("H63").Value = 2337.3333
$ws.Range("I63").Value = 2012
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 6036
$ws.Range("L63").Value = 7500
$ws.Range("M63").Value = -5287
$ws.Range("N63").Value = -8998

$ws.Range("H66").Value = 2337.3333
$ws.Range("I66").Value = 2012
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 18108
$ws.Range("L66").Value = 22500
$ws.Range("M66").Value = -14364
$ws.Range("N66").Value = -29988

$ws.Range("H107").Value = 988.6842
$ws.Range("I107").Value = 676.6667
$ws.Range("J107").Value = 1132.6923
$ws.Range("K107").Value = 2030.0001
$ws.Range("L107").Value = 3398.0769
$ws.Range("M107").Value = -110.0001
$ws.Range("N107").Value = -7238.0769

$ws.Range("H117").Value = 15699.5
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 15699.5
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 47098.5
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -53982.5

$ws.Range("H131").Value = 684844.56
$ws.Range("I131").Value = 633.63635
$ws.Range("J131").Value = 804310
$ws.Range("K131").Value = 1900.90905
$ws.Range("L131").Value = 2412930
$ws.Range("M131").Value = 3139.09095
$ws.Range("N131").Value = -2423010

$ws.Range("H132").Value = 1797.4286
$ws.Range("I132").Value = 979.8571
$ws.Range("J132").Value = 2615
$ws.Range("K132").Value = 8818.713899999999
$ws.Range("L132").Value = 23535
$ws.Range("M132").Value = -6288.713899999999
$ws.Range("N132").Value = -28595

$ws.Range("H135").Value = 1177.242
$ws.Range("I135").Value = 841.8571
$ws.Range("J135").Value = 1349.0244
$ws.Range("K135").Value = 7576.7139
$ws.Range("L135").Value = 12141.2196
$ws.Range("M135").Value = -5041.7139
$ws.Range("N135").Value = -17211.2196

$ws.Range("H139").Value = 2110.8572
$ws.Range("J139").Value = 2764.762
$ws.Range("L139").Value = 8294.286
$ws.Range("N139").Value = -18574.286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2128.4333
$ws.Range("I43").Value = 991.13043
$ws.Range("J43").Value = 5865.2856
$ws.Range("K43").Value = 991.13043
$ws.Range("L43").Value = 5865.2856
$ws.Range("M43").Value = -840.13043
$ws.Range("N43").Value = -6167.2856

$ws.Range("H46").Value = 12449.857
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 12449.857
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 12449.857
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -12761.857

$ws.Range("H80").Value = 111357330
$ws.Range("I80").Value = 200441200
$ws.Range("J80").Value = 2495
$ws.Range("K80").Value = 200441200
$ws.Range("L80").Value = 2495
$ws.Range("M80").Value = -200440202
$ws.Range("N80").Value = -4491

$ws.Range("H83").Value = 111357330
$ws.Range("I83").Value = 200441200
$ws.Range("J83").Value = 2495
$ws.Range("K83").Value = 1002206000
$ws.Range("L83").Value = 12475
$ws.Range("M83").Value = -1002201008
$ws.Range("N83").Value = -22459

$ws.Range("H126").Value = 3048.8
$ws.Range("I126").Value = 2940.923
$ws.Range("J126").Value = 3750
$ws.Range("K126").Value = 8822.769
$ws.Range("L126").Value = 11250
$ws.Range("M126").Value = -6352.769
$ws.Range("N126").Value = -16190

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4252.898
$ws.Range("I132").Value = 2396.9355
$ws.Range("J132").Value = 7449.278
$ws.Range("K132").Value = 7190.806500000001
$ws.Range("L132").Value = 22347.834
$ws.Range("M132").Value = -4660.806500000001
$ws.Range("N132").Value = -27407.834
